# aggiornamento fino a 02/05
# Appends 6 new daily rows (2021-04-27 .. 2021-05-02) to the data table,
# mirroring the date/new-cases/7-day-sum/7-day-sum-per-100k columns A:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 238 (date serial 44312 = 2021-04-26).
$lastRow = 238

# Clone column A's date formatting (style) onto the new rows first, so the
# new date cells render/format exactly like the rest of the column.
$ws.Range("A$lastRow").Copy() | Out-Null
$ws.Range("A239:A244").PasteSpecial(-4122) | Out-Null

$dates  = @(44313, 44314, 44315, 44316, 44317, 44318)
$newPos = @(7, 1, 4, 6, 3, 4)
$sum7   = @(45, 46, 42, 40, 34, 28)
$sum7_100k = @(250.3059294693514, 255.8682834575592, 233.618867504728, 222.4941595283124, 189.1200355990655, 155.7459116698187)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $lastRow + 1 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $newPos[$i]
    $ws.Cells.Item($r, 3).Value = $sum7[$i]
    $ws.Cells.Item($r, 4).Value = $sum7_100k[$i]
}
